# "Atualizacao de bases das ligas" - league fixture rows were re-ordered
# (adjacent match records swapped place within the sheet). For each pair
# of data rows, everything except the row's running id (col A), the
# division label (col C) and the match date (col D) needs to swap with
# its sibling row - i.e. column B (match id) and columns E..AB (teams,
# score, odds, P&L figures) trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get copied between paired rows (B, and E through AB).
# Column A (id), C (Div), D (Date) stay untouched.
$cols = @(2) + @(5..28)

# Snapshot original values (and whether each cell is a string) for every
# row involved in the re-pairing, BEFORE any writes happen.
$rows = @(13, 14, 41, 42, 44, 45, 50, 52, 95, 96, 100, 101, 108, 110, 118, 119, 126, 127, 148, 149, 152, 153, 156, 157, 158, 161, 162, 208, 209, 223, 224)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowData
}

# Mapping: new content of row <key> = old content of row <value>.
$mapping = @{
    13 = 14
    14 = 13
    41 = 42
    42 = 41
    44 = 45
    45 = 44
    50 = 52
    52 = 50
    95 = 96
    96 = 95
    100 = 101
    101 = 100
    108 = 110
    110 = 108
    118 = 119
    119 = 118
    126 = 127
    127 = 126
    148 = 149
    149 = 148
    152 = 153
    153 = 152
    156 = 158
    157 = 156
    158 = 157
    161 = 162
    162 = 161
    208 = 209
    209 = 208
    223 = 224
    224 = 223
}

foreach ($r in $rows) {
    $src = $mapping[$r]
    $srcData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcData[$c]
    }
}
